# Add a "2022-Q1" sheet (fund-holdings detail, same shape as the other
# quarterly sheets) positioned between "2021-Q4" and "总计", and update the
# "总计" (totals) sheet with a new leading row summarising 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet right after "2021-Q4", by copying the
#    "2021-Q4" sheet's layout/formatting (header + A-column styles, page
#    margins) and then overwriting the data with the 2022-Q1 figures.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")

$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Re-fetch "总计" *after* inserting the new sheet: worksheet references are
# resolved by position under the hood, so a handle captured before the
# insert would silently point at the new sheet once indices shift.
$total = $wb.Worksheets.Item("总计")

$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Copy the formatted header/index range (values + styles) as a starting
# point, then overwrite every data cell with the real 2022-Q1 values.
$q4.Range("A1:H3").Copy($q1.Range("A1"))

$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

$q1.Cells.Item(2,1).Value = 0
$q1.Cells.Item(2,2).NumberFormat = "@"
$q1.Cells.Item(2,2).Value = "012348"
$q1.Cells.Item(2,3).Value = "天弘恒生科技指数型发起式证券投资基金（QDII）A"
$q1.Cells.Item(2,4).NumberFormat = "@"
$q1.Cells.Item(2,4).Value = "38.10"
$q1.Cells.Item(2,5).NumberFormat = "@"
$q1.Cells.Item(2,5).Value = "92.34"
$q1.Cells.Item(2,6).NumberFormat = "@"
$q1.Cells.Item(2,6).Value = "4.72"
$q1.Cells.Item(2,7).NumberFormat = "@"
$q1.Cells.Item(2,7).Value = "1.7983"
$q1.Cells.Item(2,8).Value = 9

$q1.Cells.Item(3,1).Value = 1
$q1.Cells.Item(3,2).NumberFormat = "@"
$q1.Cells.Item(3,2).Value = "012349"
$q1.Cells.Item(3,3).Value = "天弘恒生科技指数型发起式证券投资基金（QDII）C"
$q1.Cells.Item(3,4).NumberFormat = "@"
$q1.Cells.Item(3,4).Value = "14.77"
$q1.Cells.Item(3,5).NumberFormat = "@"
$q1.Cells.Item(3,5).Value = "92.34"
$q1.Cells.Item(3,6).NumberFormat = "@"
$q1.Cells.Item(3,6).Value = "4.72"
$q1.Cells.Item(3,7).NumberFormat = "@"
$q1.Cells.Item(3,7).Value = "0.6971"
$q1.Cells.Item(3,8).Value = 9

# ---------------------------------------------------------------------
# 2. Update "总计": push the existing two rows down one row (keeping their
#    formatting) and insert the new 2022-Q1 summary row at the top.
# ---------------------------------------------------------------------
$total.Range("A2:D3").Copy($total.Range("A3:D4"))

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 2.5

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
